{"js": "// Replace the 25 \"a\u00f7b=\" division prompts in the worksheet table with new\n// values, matched positionally (table row/col) to the original document so\n// duplicate prompt text (e.g. \"19\u00f72=\" appears both as an old and a new\n// value) is never ambiguous.\nconst table = context.document.body.tables.getFirst();\n\nconst pairs = [\n  {\n    \"row\": 0,\n    \"col\": 0,\n    \"oldText\": \"63\u00f78=\",\n    \"newText\": \"89\u00f74=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 1,\n    \"oldText\": \"28\u00f74=\",\n    \"newText\": \"89\u00f74=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 2,\n    \"oldText\": \"57\u00f76=\",\n    \"newText\": \"56\u00f79=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 3,\n    \"oldText\": \"95\u00f75=\",\n    \"newText\": \"37\u00f74=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 4,\n    \"oldText\": \"36\u00f77=\",\n    \"newText\": \"88\u00f77=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 0,\n    \"oldText\": \"66\u00f72=\",\n    \"newText\": \"19\u00f72=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 1,\n    \"oldText\": \"76\u00f75=\",\n    \"newText\": \"65\u00f75=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 2,\n    \"oldText\": \"93\u00f75=\",\n    \"newText\": \"27\u00f72=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 3,\n    \"oldText\": \"84\u00f78=\",\n    \"newText\": \"23\u00f79=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 4,\n    \"oldText\": \"38\u00f78=\",\n    \"newText\": \"69\u00f73=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 0,\n    \"oldText\": \"34\u00f75=\",\n    \"newText\": \"30\u00f76=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 1,\n    \"oldText\": \"41\u00f79=\",\n    \"newText\": \"83\u00f74=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 2,\n    \"oldText\": \"76\u00f74=\",\n    \"newText\": \"86\u00f76=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 3,\n    \"oldText\": \"72\u00f77=\",\n    \"newText\": \"18\u00f79=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 4,\n    \"oldText\": \"50\u00f79=\",\n    \"newText\": \"58\u00f74=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 0,\n    \"oldText\": \"60\u00f73=\",\n    \"newText\": \"44\u00f73=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 1,\n    \"oldText\": \"15\u00f75=\",\n    \"newText\": \"14\u00f74=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 2,\n    \"oldText\": \"43\u00f75=\",\n    \"newText\": \"81\u00f73=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 3,\n    \"oldText\": \"30\u00f74=\",\n    \"newText\": \"79\u00f75=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 4,\n    \"oldText\": \"19\u00f72=\",\n    \"newText\": \"64\u00f76=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 0,\n    \"oldText\": \"12\u00f76=\",\n    \"newText\": \"17\u00f75=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 1,\n    \"oldText\": \"33\u00f75=\",\n    \"newText\": \"92\u00f73=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 2,\n    \"oldText\": \"18\u00f76=\",\n    \"newText\": \"30\u00f73=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 3,\n    \"oldText\": \"70\u00f78=\",\n    \"newText\": \"34\u00f72=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 4,\n    \"oldText\": \"26\u00f79=\",\n    \"newText\": \"72\u00f74=\"\n  }\n];\n\n// Resolve every target cell and load its current text in one batch.\nconst cells = pairs.map(({ row, col }) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\n// Verify the existing content matches what we expect before overwriting, so\n// a mismatched document fails loudly instead of silently mis-editing it.\npairs.forEach(({ row, col, oldText }, i) => {\n  const current = (cells[i].value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Unexpected cell text at [${row},${col}]: expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n});\n\n// Apply every replacement, then flush the writes in a single sync.\npairs.forEach(({ newText }, i) => {\n  cells[i].value = newText;\n});\nawait context.sync();\n", "ps1": "# Replace the 25 \"a\u00f7b=\" division prompts in the worksheet table with new\n# values, matched positionally (table row/col) to the original document so\n# duplicate prompt text (e.g. \"19\u00f72=\" appears both as an old and a new\n# value) is never ambiguous.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairs = @(\n  @{ Row = 1; Col = 1; Old = \"63\u00f78=\"; New = \"89\u00f74=\" },\n  @{ Row = 1; Col = 2; Old = \"28\u00f74=\"; New = \"89\u00f74=\" },\n  @{ Row = 1; Col = 3; Old = \"57\u00f76=\"; New = \"56\u00f79=\" },\n  @{ Row = 1; Col = 4; Old = \"95\u00f75=\"; New = \"37\u00f74=\" },\n  @{ Row = 1; Col = 5; Old = \"36\u00f77=\"; New = \"88\u00f77=\" },\n  @{ Row = 5; Col = 1; Old = \"66\u00f72=\"; New = \"19\u00f72=\" },\n  @{ Row = 5; Col = 2; Old = \"76\u00f75=\"; New = \"65\u00f75=\" },\n  @{ Row = 5; Col = 3; Old = \"93\u00f75=\"; New = \"27\u00f72=\" },\n  @{ Row = 5; Col = 4; Old = \"84\u00f78=\"; New = \"23\u00f79=\" },\n  @{ Row = 5; Col = 5; Old = \"38\u00f78=\"; New = \"69\u00f73=\" },\n  @{ Row = 9; Col = 1; Old = \"34\u00f75=\"; New = \"30\u00f76=\" },\n  @{ Row = 9; Col = 2; Old = \"41\u00f79=\"; New = \"83\u00f74=\" },\n  @{ Row = 9; Col = 3; Old = \"76\u00f74=\"; New = \"86\u00f76=\" },\n  @{ Row = 9; Col = 4; Old = \"72\u00f77=\"; New = \"18\u00f79=\" },\n  @{ Row = 9; Col = 5; Old = \"50\u00f79=\"; New = \"58\u00f74=\" },\n  @{ Row = 13; Col = 1; Old = \"60\u00f73=\"; New = \"44\u00f73=\" },\n  @{ Row = 13; Col = 2; Old = \"15\u00f75=\"; New = \"14\u00f74=\" },\n  @{ Row = 13; Col = 3; Old = \"43\u00f75=\"; New = \"81\u00f73=\" },\n  @{ Row = 13; Col = 4; Old = \"30\u00f74=\"; New = \"79\u00f75=\" },\n  @{ Row = 13; Col = 5; Old = \"19\u00f72=\"; New = \"64\u00f76=\" },\n  @{ Row = 17; Col = 1; Old = \"12\u00f76=\"; New = \"17\u00f75=\" },\n  @{ Row = 17; Col = 2; Old = \"33\u00f75=\"; New = \"92\u00f73=\" },\n  @{ Row = 17; Col = 3; Old = \"18\u00f76=\"; New = \"30\u00f73=\" },\n  @{ Row = 17; Col = 4; Old = \"70\u00f78=\"; New = \"34\u00f72=\" },\n  @{ Row = 17; Col = 5; Old = \"26\u00f79=\"; New = \"72\u00f74=\" }\n)\n\nforeach ($p in $pairs) {\n  $cell = $t.Cell($p.Row, $p.Col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $p.Old) {\n    throw \"Unexpected cell text at [$($p.Row),$($p.Col)]: expected '$($p.Old)' but found '$current'\"\n  }\n  $cell.Range.Text = $p.New\n}\n"}
